$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.103.85"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +4.66%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.248.65"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +2.38%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "396.12"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.26%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.16"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.60%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.590"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +7.38%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.244.07"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +2.41%  "

$ws.Range("E9").Value = "  +0.03%  "

$ws.Range("E10").Value = "  +1.44%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "39.26"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.28%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0986"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +9.54%  "

$ws.Range("E13").Value = "  +2.12%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.757.51"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +2.26%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.36"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +3.58%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.09"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.04%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.253.83"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +2.50%  "

$ws.Range("E18").Value = "  -2.77%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.76"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.90%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "56.948.15"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +4.61%  "

$ws.Range("E21").Value = "  +1.38%  "

$ws.Range("E22").Value = "  +7.75%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.08"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.95%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "294.91"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +7.05%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "74.22"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +2.88%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.17"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.62%  "

$ws.Range("E27").Value = "  +0.89%  "

$ws.Range("E28").Value = "  +0.98%  "

$ws.Range("E29").Value = "  -5.36%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.22"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -5.59%  "

$ws.Range("E31").Value = "  -1.11%  "

$ws.Range("E32").Value = "  -0.04%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.23"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.92%  "

$ws.Range("E34").Value = "  -3.12%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "39.86"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +8.79%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0488"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -3.58%  "

$ws.Range("E37").Value = "  +1.38%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "51.55"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.65%  "

$ws.Range("E39").Value = "  -0.13%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.47"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -4.57%  "

$ws.Range("E41").Value = "  +2.01%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "137.94"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +4.27%  "

$ws.Range("E43").Value = "  +3.99%  "

$ws.Range("E44").Value = "  -1.82%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "17.08"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.08%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.94"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -3.57%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.280"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -4.67%  "

$ws.Range("B48").Value = "WEMIXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.32"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +11.53%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "22.21"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.41%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.162.19"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +2.83%  "

$ws.Range("E51").Value = "  -6.12%  "
